$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy header style (H1, bold/centered/bordered) onto the new I1:J1 header cells
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Fill in the I0 / IF numeric data for rows 2-73
$iVals = @(8,8,9,7,9,9,8,10,9,9,9,9,8,9,9,9,9,9,9,9,9,8,8,9,6,9,9,9,9,9,7,9,9,9,9,9,9,9,8,9,8,8,7,8,8,5,9,9,8,8,9,10,9,9,8,6,8,9,8,8,9,9,8,8,7,7,2,9,8,4,8,5)
$jVals = @(8,9,9,7,9,9,8,10,9,9,9,9,9,9,9,9,9,9,9,9,9,9,9,9,6,9,9,9,9,9,7,9,9,9,9,9,9,10,8,9,8,8,8,9,8,5,9,9,8,9,9,10,9,9,8,6,8,9,8,8,9,9,8,8,7,7,3,9,8,4,8,5)
for ($r = 2; $r -le 73; $r++) {
    $ws.Cells.Item($r, 9).Value = $iVals[$r - 2]
    $ws.Cells.Item($r, 10).Value = $jVals[$r - 2]
}

Write-Output "done"